# SA Attendance.xlsx - fill in "Session 17" (column W) attendance for rows 7-60.
#
# Column W was previously blank for every participant (style 41, no value).
# This edit records their Session 17 attendance as "P" (Present) or "A"
# (Absent), matching the same "P"/"A" text + style used by the other
# session columns (e.g. column V, style 40). Row 50's participant joined
# partway through the term (no attendance formulas / earlier sessions
# recorded), so only their W50 cell gets a value. Row 49 additionally had
# an earlier entry (N49) corrected from "A" to "P".
#
# The Total Absence (column E) / Total Present (column F) COUNTIF formulas
# recalculate automatically once the underlying P/A cells change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rows where Session 17 attendance = Present ("P")
$presentRows = @(
  7,8,9,10,11,12,13,14,15,16,
  18,20,21,22,23,24,25,
  27,29,
  31,32,34,
  36,37,39,
  41,43,44,45,46,48,
  49,50,51,52,53,54,55,56,57,58,59,60
)

# Rows where Session 17 attendance = Absent ("A")
$absentRows = @(17,19,26,28,30,33,35,38,40,42,47)

# Copy the formatting already used for filled-in attendance cells (e.g. V7,
# style 40) onto the whole W7:W60 block, so the newly-entered cells pick up
# the correct (non-blank) cell style instead of keeping the blank style 41.
$ws.Range("V7").Copy()
$ws.Range("W7:W60").PasteSpecial(-4122)  # xlPasteFormats

foreach ($r in $presentRows) {
  $ws.Range("W$r").Value = "P"
}
foreach ($r in $absentRows) {
  $ws.Range("W$r").Value = "A"
}

# Row 49 also corrects an earlier recorded session (column N) from Absent to Present.
$ws.Range("N49").Value = "P"

Write-Host "Filled Session 17 (column W) attendance for rows 7-60 and corrected N49."
